$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header for column C ---
$ws.Range("C1").Value = "% Crecimiento"
$ws.Range("C1").Font.Bold = $true

# --- Growth-rate formulas, row 3 through row 21 ---
# (B_n / B_(n-1) - 1) / (A_n - A_(n-1))
for ($r = 3; $r -le 21; $r++) {
    $prev = $r - 1
    $ws.Range("C$r").Formula = "=(B$r/B$prev-1)/(A$r-A$prev)"
}

# Percentage formatting (0.00%) for the growth column
$ws.Range("C3:C21").Style = "Percent"
$ws.Range("C3:C21").NumberFormat = "0.00%"

# --- Average of the last 30 years worth of growth (rows 15-21) ---
$ws.Range("B23").Value = "Promedio últimos 30 años"
$ws.Range("B23").Font.Bold = $true

$ws.Range("C23").Formula = "=AVERAGE(C15:C21)"
$ws.Range("C23").Style = "Percent"
$ws.Range("C23").NumberFormat = "0.00%"
$ws.Range("C23").Font.Bold = $true

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 24.7109375
$ws.Columns.Item(3).ColumnWidth = 14

# --- Selection, matching the authored workbook state ---
$ws.Range("C23").Select()
